# Update the "Metadata" sheet (Property/Value pairs).
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/offset-end"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet: the root "Extension" row (row 2) should no longer
# carry the ele-1/ext-1 constraint text in its Constraint(s) column (AI) - that
# text now only belongs to the "Extension.extension" row (row 4). Copy a
# neighboring cell that already holds an (empty) text string into AI2 so the
# cell keeps its text type instead of collapsing into a truly blank cell.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI3").Copy($elements.Range("AI2"))
